$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2 through 18 from 2023-09-01 to 2023-09-05
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = Get-Date -Year 2023 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
}
